# Applies the scheduled-runner market-data refresh to the Sheets workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 549.2381
$ws.Range("J17").Value = 546.46344
$ws.Range("L17").Value = 1639.39032
$ws.Range("N17").Value = -1975.39032

$ws.Range("H40").Value = 1222.0625
$ws.Range("I40").Value = 950.3333
$ws.Range("J40").Value = 1284.7693
$ws.Range("K40").Value = 950.3333
$ws.Range("L40").Value = 1284.7693
$ws.Range("M40").Value = -775.3333
$ws.Range("N40").Value = -1634.7693

$ws.Range("H53").Value = 3140
$ws.Range("I53").Value = 383
$ws.Range("K53").Value = 383
$ws.Range("M53").Value = 254

$ws.Range("H76").Value = 3706809.8
$ws.Range("I76").Value = 3237.5
$ws.Range("J76").Value = 7939463.5
$ws.Range("K76").Value = 3237.5
$ws.Range("L76").Value = 7939463.5
$ws.Range("M76").Value = -2922.5
$ws.Range("N76").Value = -7940093.5

$ws.Range("H79").Value = 3706809.8
$ws.Range("I79").Value = 3237.5
$ws.Range("J79").Value = 7939463.5
$ws.Range("K79").Value = 3237.5
$ws.Range("L79").Value = 7939463.5
$ws.Range("M79").Value = -2145.5
$ws.Range("N79").Value = -7941647.5

$ws.Range("H129").Value = 964.90247
$ws.Range("I129").Value = 532.3333
$ws.Range("J129").Value = 999.0526
$ws.Range("K129").Value = 1596.9999
$ws.Range("L129").Value = 2997.1578
$ws.Range("M129").Value = 3403.0001
$ws.Range("N129").Value = -12997.1578

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3175.7856
$ws.Range("I32").Value = 1810.5428
$ws.Range("K32").Value = 1810.5428
$ws.Range("M32").Value = -1523.5428

$ws.Range("H61").Value = 7681.5
$ws.Range("I61").Value = 11336
$ws.Range("J61").Value = 5071.143
$ws.Range("K61").Value = 11336
$ws.Range("L61").Value = 5071.143
$ws.Range("M61").Value = -11124
$ws.Range("N61").Value = -5495.143

$ws.Range("H88").Value = 61850.176
$ws.Range("I88").Value = 2001.6666
$ws.Range("J88").Value = 74674.86
$ws.Range("K88").Value = 2001.6666
$ws.Range("L88").Value = 74674.86
$ws.Range("M88").Value = -1595.6666
$ws.Range("N88").Value = -75486.86

$ws.Range("H91").Value = 61850.176
$ws.Range("I91").Value = 2001.6666
$ws.Range("J91").Value = 74674.86
$ws.Range("K91").Value = 2001.6666
$ws.Range("L91").Value = 74674.86
$ws.Range("M91").Value = -597.6666
$ws.Range("N91").Value = -77482.86

$ws.Range("H132").Value = 15549.583
$ws.Range("I132").Value = 1131.4
$ws.Range("K132").Value = 3394.2
$ws.Range("M132").Value = -864.2000000000003

$ws.Range("H136").Value = 7681.5
$ws.Range("I136").Value = 11336
$ws.Range("J136").Value = 5071.143
$ws.Range("K136").Value = 34008
$ws.Range("L136").Value = 15213.429
$ws.Range("M136").Value = -31458
$ws.Range("N136").Value = -20313.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 266.07693
$ws.Range("I22").Value = 266.07693
$ws.Range("K22").Value = 266.07693
$ws.Range("M22").Value = -93.07693

$ws.Range("H86").Value = 1489.2222
$ws.Range("I86").Value = 1233.3334
$ws.Range("J86").Value = 2001
$ws.Range("K86").Value = 1233.3334
$ws.Range("L86").Value = 2001
$ws.Range("M86").Value = -110.3334
$ws.Range("N86").Value = -4247

$ws.Range("H89").Value = 1489.2222
$ws.Range("I89").Value = 1233.3334
$ws.Range("J89").Value = 2001
$ws.Range("K89").Value = 6166.666999999999
$ws.Range("L89").Value = 10005
$ws.Range("M89").Value = -550.6669999999995
$ws.Range("N89").Value = -21237

$ws.Range("H134").Value = 2811.647
$ws.Range("I134").Value = 2986.5334
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 8959.600199999999
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -6424.600199999999
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 121.125
$ws.Range("J7").Value = 182.5
$ws.Range("L7").Value = 182.5
$ws.Range("N7").Value = -408.5

$ws.Range("H31").Value = 2262.6758
$ws.Range("I31").Value = 1008.36
$ws.Range("J31").Value = 4875.8335
$ws.Range("K31").Value = 1008.36
$ws.Range("L31").Value = 4875.8335
$ws.Range("M31").Value = -713.36
$ws.Range("N31").Value = -5465.8335

$ws.Range("H34").Value = 2262.6758
$ws.Range("I34").Value = 1008.36
$ws.Range("J34").Value = 4875.8335
$ws.Range("K34").Value = 1008.36
$ws.Range("L34").Value = 4875.8335
$ws.Range("M34").Value = -806.36
$ws.Range("N34").Value = -5279.8335

$ws.Range("H62").Value = 4539.8
$ws.Range("I62").Value = 4016.5
$ws.Range("J62").Value = 5324.75
$ws.Range("K62").Value = 4016.5
$ws.Range("L62").Value = 5324.75
$ws.Range("M62").Value = -3392.5
$ws.Range("N62").Value = -6572.75

$ws.Range("H65").Value = 4539.8
$ws.Range("I65").Value = 4016.5
$ws.Range("J65").Value = 5324.75
$ws.Range("K65").Value = 20082.5
$ws.Range("L65").Value = 26623.75
$ws.Range("M65").Value = -16962.5
$ws.Range("N65").Value = -32863.75

$ws.Range("H99").Value = 5142.2383
$ws.Range("I99").Value = 3829.818
$ws.Range("J99").Value = 6585.9
$ws.Range("K99").Value = 3829.818
$ws.Range("L99").Value = 6585.9
$ws.Range("M99").Value = -2331.818
$ws.Range("N99").Value = -9581.9

$ws.Range("H126").Value = 5142.2383
$ws.Range("I126").Value = 3829.818
$ws.Range("J126").Value = 6585.9
$ws.Range("K126").Value = 11489.454
$ws.Range("L126").Value = 19757.7
$ws.Range("M126").Value = -9019.454000000002
$ws.Range("N126").Value = -24697.7

$ws.Range("H132").Value = 3082.7273
$ws.Range("I132").Value = 899.1667
$ws.Range("J132").Value = 5703
$ws.Range("K132").Value = 2697.5001
$ws.Range("L132").Value = 17109
$ws.Range("M132").Value = -167.5001000000002
$ws.Range("N132").Value = -22169

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1162.742
$ws.Range("J68").Value = 1168.1666
$ws.Range("L68").Value = 3504.4998
$ws.Range("N68").Value = -5126.4998

$ws.Range("H71").Value = 1162.742
$ws.Range("J71").Value = 1168.1666
$ws.Range("L71").Value = 10513.4994
$ws.Range("N71").Value = -18625.4994

$ws.Range("H117").Value = 1373
$ws.Range("I117").Value = 661.75
$ws.Range("J117").Value = 2321.3333
$ws.Range("K117").Value = 1985.25
$ws.Range("L117").Value = 6963.999899999999
$ws.Range("M117").Value = 1456.75
$ws.Range("N117").Value = -13847.9999

$ws.Range("H134").Value = 5779.522
$ws.Range("I134").Value = 5901.619
$ws.Range("J134").Value = 4497.5
$ws.Range("K134").Value = 17704.857
$ws.Range("L134").Value = 13492.5
$ws.Range("M134").Value = -12634.857
$ws.Range("N134").Value = -23632.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6948499
$ws.Range("I70").Value = 4300
$ws.Range("K70").Value = 4300
$ws.Range("M70").Value = -4030

$ws.Range("H73").Value = 6948499
$ws.Range("I73").Value = 4300
$ws.Range("K73").Value = 4300
$ws.Range("M73").Value = -3364

$ws.Range("H132").Value = 43735
$ws.Range("I132").Value = 5903
$ws.Range("J132").Value = 104266.2
$ws.Range("K132").Value = 17709
$ws.Range("L132").Value = 312798.6
$ws.Range("M132").Value = -15179
$ws.Range("N132").Value = -317858.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1197.5
$ws.Range("I22").Value = 1195
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 1195
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -900
$ws.Range("N22").Value = -1790

$ws.Range("H27").Value = 1197.5
$ws.Range("I27").Value = 1195
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 1195
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = -1088
$ws.Range("N27").Value = -1414

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
